# Mifos -> Finflux: insert a new blank column before column N ("Late")
# on the "Repayment schedule" sheet, shifting the existing Late / heading /
# Outstanding columns one place to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new column at N; existing N/O/P data shifts to O/P/Q.
$ws.Range("N1").EntireColumn.Insert()

# Match the width of the newly inserted (blank) column to its left
# neighbour, mirroring how Excel carries formatting into inserted columns.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Restore the selection to reflect where the user ended up after editing.
$ws.Range("Q5").Select()
